$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ---
$ws.Range("A5").Value = 42395
$ws.Range("A5").NumberFormat = "d-mmm"

$ws.Range("B5").Value = 0.73055555555555562
$ws.Range("B5").NumberFormat = "h:mm"

$ws.Range("C5").Value = 0.75694444444444453
$ws.Range("C5").NumberFormat = "h:mm"

$ws.Range("D5").Formula = "=C5-B5"
$ws.Range("D5").NumberFormat = "h:mm"

$ws.Range("E5").Value = "ListView, Database Connection"

# --- Row 6 ---
# Reuse the existing date style (same as A2:A4) by copying it over,
# rather than re-setting NumberFormat (which would mint a duplicate
# custom numFmt instead of reusing the built-in one).
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 42395

$ws.Range("B6").Value = 0.77430555555555547
$ws.Range("B6").NumberFormat = "h:mm"

$ws.Range("C6").Value = 0.80347222222222225
$ws.Range("C6").NumberFormat = "h:mm"

$ws.Range("D6").Formula = "=C6-B6"
$ws.Range("D6").NumberFormat = "h:mm"

$ws.Range("E6").Value = "Database Connection"

# --- Selection update ---
$ws.Range("E9").Select()
